# Final commit from Jeff for this iteration
# Updates the Iteration #3 "Work Effort" column (E29:E43) with actual hours,
# and clears out the placeholder "X" marks previously in column H,
# and updates the view selection/scroll position.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Column E: fill in work-effort hours for rows 29-43 ---
$ws.Range("E29").Value = 4
$ws.Range("E30").Value = 6
$ws.Range("E31").Value = 1.5
$ws.Range("E32").Value = 2.5
$ws.Range("E33").Value = 2
$ws.Range("E34").Value = 2
$ws.Range("E35").Value = 2
$ws.Range("E36").Value = 4
$ws.Range("E37").Value = 3
$ws.Range("E38").Value = 2
$ws.Range("E39").Value = 2
$ws.Range("E41").Value = 32
$ws.Range("E42").Value = 2
$ws.Range("E43").Value = 3

# --- Column H: clear the old "X" placeholders / leftover values (rows 29-43) ---
# Rows with a percentage-formatted H cell (style 17) keep the (now empty) cell
# so its number format is preserved; the rest had the "X" string with no
# distinguishing format, so the whole cell is cleared away.
$ws.Range("H29").Clear()
$ws.Range("H30").ClearContents()
$ws.Range("H31").Clear()
$ws.Range("H32").Clear()
$ws.Range("H33").Clear()
$ws.Range("H34").Clear()
$ws.Range("H35").ClearContents()
$ws.Range("H36").ClearContents()
$ws.Range("H37").Clear()
$ws.Range("H38").Clear()
$ws.Range("H39").Clear()
$ws.Range("H40").Clear()
$ws.Range("H41").Clear()
$ws.Range("H42").Clear()
$ws.Range("H43").ClearContents()

# --- Update the saved view state (scroll position / selection) ---
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 16
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("F40").Select()
